# Logging for debugging purposes
Write-Output "Starting schedule update script"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Output "Loaded active workbook and worksheet"

# New week's schedule data: Teams, Time, Location
$schedule = @(
    @("Minnesota  @  LA Rams", "8:15 PM", "SoFi Stadium"),
    @("Green Bay  @  Jacksonville", "1:00 PM", "EverBank Stadium"),
    @("Atlanta  @  Tampa Bay", "1:00 PM", "Raymond James Stadium"),
    @("NY Jets  @  New England", "1:00 PM", "Gillette Stadium"),
    @("Arizona  @  Miami", "1:00 PM", "Hard Rock Stadium"),
    @("Indianapolis  @  Houston", "1:00 PM", "NRG Stadium"),
    @("Tennessee  @  Detroit", "1:00 PM", "Ford Field"),
    @("Baltimore  @  Cleveland", "1:00 PM", "Huntington Bank Field"),
    @("Philadelphia  @  Cincinnati", "1:00 PM", "Paycor Stadium"),
    @("New Orleans  @  LA Chargers", "4:05 PM", "SoFi Stadium"),
    @("Buffalo  @  Seattle", "4:05 PM", "Lumen Field"),
    @("Kansas City  @  Las Vegas", "4:25 PM", "Allegiant Stadium"),
    @("Carolina  @  Denver", "4:25 PM", "Empower Field at Mile High"),
    @("Chicago  @  Washington", "4:25 PM", "Northwest Stadium"),
    @("Dallas  @  San Francisco", "8:20 PM", "Levi's Stadium"),
    @("NY Giants  @  Pittsburgh", "8:15 PM", "Acrisure Stadium")
)

Write-Output ("Schedule contains " + $schedule.Length + " games")

$row = 2
foreach ($game in $schedule) {
    Write-Output ("Writing row " + $row + ": " + $game[0])
    $ws.Cells.Item($row, 1).Value = $game[0]
    $ws.Cells.Item($row, 2).Value = $game[1]
    $ws.Cells.Item($row, 3).Value = $game[2]
    $row = $row + 1
}

Write-Output "Finished writing schedule rows"

Write-Output "Schedule update script complete"
